$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (record #22) was an empty placeholder row; fill it in as a new
# expense entry: 支出 / 300 / 2017-11-21 / 生活费 / 生活费(12/01-12/10).
# Copy formatting from row 22 (an existing, fully-populated expense row)
# so the new row picks up the correct fill/border/alignment styles.
$ws.Range("B22:G22").Copy()
$ws.Range("B24:G24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C24").Value = "支出"
$ws.Range("D24").Value = 300
$ws.Range("E24").Value = "11/21/2017"
$ws.Range("F24").Value = "生活费"
$ws.Range("G24").Value = "生活费(12/01-12/10)"

# The saved view's active cell moved to G26 after the edit.
$ws.Range("G26").Select()
